# Add a new "Save" column (H) to the s_vals sheet, matching the
# existing header formatting used by the other stat columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (G1) onto the
# new header cell H1 so the new column reuses the same cell style
# (bold, centered, top-aligned, thin border) instead of creating a
# brand new style entry.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new column's values.
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
